$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add resolution notes in column C for several defects (bug triage pass of 2/2/2020)
$ws.Range("C15").Value = "[2/2/2020]Done"
$ws.Range("C17").Value = "[2/2/2020]Could not recreate"
$ws.Range("C18").Value = "[2/2/2020]Fixed"
$ws.Range("C19").Value = "[2/2/2020]Added hardcoded values"
$ws.Range("C21").Value = "[2/2/2020]Fixed"

# Narrow column C now that the new notes are shorter than the old ones
$ws.Columns.Item(3).ColumnWidth = 55.33

# Leave the active cell on B19, matching where the user ended up after editing
$ws.Range("B19").Select()
